# Update de bases das ligas, do dia 24-02-2024 as 23:13
# Swaps the data of two match rows that had their order corrected
# (columns B:AC, keeping column A - the sequence number - untouched),
# plus a small odds correction on row 193 (columns U/V).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Swap-Rows($sheet, $row1, $row2, $columns) {
    foreach ($col in $columns) {
        $cell1 = $sheet.Range("$col$row1")
        $cell2 = $sheet.Range("$col$row2")
        $val1 = $cell1.Value2
        $val2 = $cell2.Value2
        $cell1.Value2 = $val2
        $cell2.Value2 = $val1
    }
}

# Rows 161 and 162 had their match data swapped
Swap-Rows $ws 161 162 $cols

# Rows 168 and 169 had their match data swapped
Swap-Rows $ws 168 169 $cols

# Row 193: corrected closing odds for oddAHOver / oddAHUnder
$ws.Range("U193").Value2 = 1.8
$ws.Range("V193").Value2 = 2
